$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the existing row 239 ("Clementina"/45119 row),
# shifting the old rows 239-247 down to 241-249.
$ws.Rows.Item(239).Insert()
$ws.Rows.Item(239).Insert()

# --- New row 239 ---
$ws.Range("A239").Value = 11
$ws.Range("B239").Value = "Vega Monumental Concepción"
$ws.Range("C239").Value = "Bíobío"
$ws.Range("D239").Value = 45147
$ws.Range("E239").Value = 8
$ws.Range("F239").Value = "Fruta"
$ws.Range("G239").Value = 100102
$ws.Range("H239").Value = "Cítricos"
$ws.Range("I239").Value = 100102004
$ws.Range("J239").Value = "Mandarina"
$ws.Range("K239").Value = "Clemenuless"
$ws.Range("L239").Value = "Primera"
$ws.Range("M239").Value = 250
$ws.Range("N239").Value = 7500
$ws.Range("O239").Value = 8000
$ws.Range("P239").Value = 7800
$ws.Range("Q239").Value = "$/bandeja 10 kilos"
$ws.Range("R239").Value = "Región de O'Higgins"
$ws.Range("S239").Value = 780
$ws.Range("T239").Value = 10

# --- New row 240 ---
$ws.Range("A240").Value = 11
$ws.Range("B240").Value = "Vega Monumental Concepción"
$ws.Range("C240").Value = "Bíobío"
$ws.Range("D240").Value = 45147
$ws.Range("E240").Value = 8
$ws.Range("F240").Value = "Fruta"
$ws.Range("G240").Value = 100102
$ws.Range("H240").Value = "Cítricos"
$ws.Range("I240").Value = 100102004
$ws.Range("J240").Value = "Mandarina"
$ws.Range("K240").Value = "Murcott"
$ws.Range("L240").Value = "Primera"
$ws.Range("M240").Value = 220
$ws.Range("N240").Value = 10000
$ws.Range("O240").Value = 11000
$ws.Range("P240").Value = 10545
$ws.Range("Q240").Value = "$/bandeja 18 kilos"
$ws.Range("R240").Value = "Región de O'Higgins"
$ws.Range("S240").Value = 586
$ws.Range("T240").Value = 18
